# ---------------------------------------------------------------------------
# Update integration tests with achievements
# ---------------------------------------------------------------------------
$d = $word.ActiveDocument

$LDQ = [char]8220   # "  left double quote
$RDQ = [char]8221   # "  right double quote

# ---------------------------------------------------------------------------
# 1. Collapse "R" + "eturn" + "ed" + " their profile page URL..." into a
#    single run.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    ("R" + "eturn" + "ed" + " their profile page URL and response of 200"),
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Returned their profile page URL and response of 200", 2)

# ---------------------------------------------------------------------------
# 2. Collapse "<RDQ> return" + "ed" + " <LDQ>True<RDQ>" into a single run.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    ($RDQ + " return" + "ed" + " " + $LDQ + "True" + $RDQ),
    $false, $false, $false, $false, $false, $true, 1, $false,
    ($RDQ + " returned " + $LDQ + "True" + $RDQ), 2)

# ---------------------------------------------------------------------------
# 3. Collapse "Inc" + "orrect image file type is trying to be stored".
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    ("Inc" + "orrect image file type is trying to be stored"),
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Incorrect image file type is trying to be stored", 2)

# ---------------------------------------------------------------------------
# 4. Collapse "<RDQ> returns <LDQ>" + "False" + "<RDQ>" into a single run.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    ($RDQ + " returns " + $LDQ + "False" + $RDQ),
    $false, $false, $false, $false, $false, $true, 1, $false,
    ($RDQ + " returns " + $LDQ + "False" + $RDQ), 2)

# ---------------------------------------------------------------------------
# 5. Resize the "Integration Tests" table's 2nd and 4th columns (widen the
#    comments column, narrow the components column) - this is table index 4
#    (1-based) / Tables.Item(4).
# ---------------------------------------------------------------------------
$integrationTable = $d.Tables.Item(4)

$integrationTable.Columns.Item(2).Width = 174.8
$integrationTable.Columns.Item(4).Width = 189.95

$CR = [char]13
$LSQ = [char]8216   # '  left single quote
$RSQ = [char]8217   # '  right single quote

function Add-BlankRow($table) {
    $row = $table.Rows.Add()
    return $row
}

function Set-SuperscriptAfterFind($cell, [string]$needle, [int]$skipChars) {
    # Finds `needle` inside `cell`, then makes everything from `skipChars`
    # characters into the match through the end of the match superscript
    # (e.g. needle="10th", skipChars=2 -> only "th" becomes superscript).
    # Rebuilding the range via $d.Range(start,end) (rather than reusing the
    # Cell's own .Range object across Find calls) avoids a stale-match bug
    # in this host where a 2nd Find on a different cell's range still
    # reports the previous cell's match offsets.
    $s = $cell.Range.Start
    $e = $cell.Range.End
    $r = $d.Range($s, $e)
    $found = $r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $sub = $d.Range($r.Start + $skipChars, $r.End)
        $sub.Font.Superscript = $true
    }
    return $found
}

# ---------------------------------------------------------------------------
# 6. Blank spacer row immediately after row 14.
# ---------------------------------------------------------------------------
$null = Add-BlankRow $integrationTable

# ---------------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------------
$row15 = Add-BlankRow $integrationTable
$row15.Cells.Item(1).Range.Text = "15"
$row15.Cells.Item(2).Range.Text = "Viewing your profile earns achievement ID 1"
$row15.Cells.Item(3).Range.Text = "Success"
$row15.Cells.Item(4).Range.Text = (
    "The username and achievement ID were added to the CompleteAchievements table and 25xp was added for the user in the UserLevel table, as expected. Nothing was awarded on subsequent profile views." +
    $CR +
    "Achievement ID reference table can be found as " + $LSQ + "Achievements_table_reference.csv" + $RSQ + "."
)

# ---------------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------------
$row16 = Add-BlankRow $integrationTable
$row16.Cells.Item(1).Range.Text = "16"
$row16.Cells.Item(2).Range.Text = "Adding a connection for the first time completes achievement ID 4 for user and the user you are connecting to if necessary"
$row16.Cells.Item(3).Range.Text = "Success"
$row16.Cells.Item(4).Range.Text = "Since neither the connection invitee nor recipient had the achievement, both received it in the database as well as the relevant xp."

# ---------------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------------
$row17 = Add-BlankRow $integrationTable
$row17.Cells.Item(1).Range.Text = "17"
$row17.Cells.Item(2).Range.Text = "Adding a connection completes achievement ID 5 for user or the user you are connecting to if it is their 10th connection"
$row17.Cells.Item(3).Range.Text = "Success"
$row17.Cells.Item(4).Range.Text = "Tested with the requestee and the requester gaining their 10th connection separately and in the same request acceptance."
$null = Set-SuperscriptAfterFind $row17.Cells.Item(2) "10th" 2
$null = Set-SuperscriptAfterFind $row17.Cells.Item(4) "10th" 2

# ---------------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------------
$row18 = Add-BlankRow $integrationTable
$row18.Cells.Item(1).Range.Text = "18"
$row18.Cells.Item(2).Range.Text = "Posting to the social feed for the first time completes achievement ID 7"
$row18.Cells.Item(3).Range.Text = "Success"
$row18.Cells.Item(4).Range.Text = "No achievement for a second post"

# ---------------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------------
$row19 = Add-BlankRow $integrationTable
$row19.Cells.Item(1).Range.Text = "19"
$row19.Cells.Item(2).Range.Text = "Posting to the social feed for the 20th time completes achievement ID 9"
$row19.Cells.Item(3).Range.Text = "Success"
$null = Set-SuperscriptAfterFind $row19.Cells.Item(2) "20th" 2

# ---------------------------------------------------------------------------
# Row 20
# ---------------------------------------------------------------------------
$row20 = Add-BlankRow $integrationTable
$row20.Cells.Item(1).Range.Text = "20"
$row20.Cells.Item(2).Range.Text = "Commenting on a post for the first time completes achievement ID 10"
$row20.Cells.Item(3).Range.Text = "Success"
$row20.Cells.Item(4).Range.Text = "No achievement for subsequent comments"

# ---------------------------------------------------------------------------
# Row 21
# ---------------------------------------------------------------------------
$row21 = Add-BlankRow $integrationTable
$row21.Cells.Item(1).Range.Text = "21"
$row21.Cells.Item(2).Range.Text = "Adding a close connection for the first time completes achievement ID 12 for user"
$row21.Cells.Item(3).Range.Text = "Success"
$row21.Cells.Item(4).Range.Text = "Close friends is a one way system, user1 marks user2 as user1" + $RSQ + "s close friend only. Therefore, this achievement is not awarded to user 2 in this case "

# ---------------------------------------------------------------------------
# Row 22
# ---------------------------------------------------------------------------
$row22 = Add-BlankRow $integrationTable
$row22.Cells.Item(1).Range.Text = "22"
$row22.Cells.Item(2).Range.Text = "Adding a connection with a mutual interest on their profile completes achievement ID 16 for user or the user you are connecting to if it is their first of such connections"
$row22.Cells.Item(3).Range.Text = "Success"
$row22.Cells.Item(4).Range.Text = "Both the connection requestee and the requester gained this achievement"

# ---------------------------------------------------------------------------
# Row 23
# ---------------------------------------------------------------------------
$row23 = Add-BlankRow $integrationTable
$row23.Cells.Item(1).Range.Text = "23"
$row23.Cells.Item(2).Range.Text = "Opening a post and pressing like for the first time completes achievement ID 19 for user"
$row23.Cells.Item(3).Range.Text = "Success"
$row23.Cells.Item(4).Range.Text = (
    "The user was also able to unlike the post; however this did not allow them to get the achievement again." +
    $CR +
    "The user who posted the liked post is also checked foe needing achievement IDs 20 and 22 relating to receiving likes."
)

# ---------------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------------
$row24 = Add-BlankRow $integrationTable
$row24.Cells.Item(1).Range.Text = "24"
$row24.Cells.Item(2).Range.Text = "Completing a quiz for the first time as well as with achieving a perfect score should award achievement IDs 27 and 28 to user"
$row24.Cells.Item(3).Range.Text = "Success"
$row24.Cells.Item(4).Range.Text = "On completing quiz, achievement ID 27 is awarded, and since the score is 100%, ID 28 is awarded also. Achievement ID 28 can also be attained on a later quiz attempt "

# ---------------------------------------------------------------------------
# Two trailing blank spacer rows.
# ---------------------------------------------------------------------------
$null = Add-BlankRow $integrationTable
$null = Add-BlankRow $integrationTable
